$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 6050
$ws.Range("I51").Value = 2000
$ws.Range("J51").Value = 7400
$ws.Range("K51").Value = 2000
$ws.Range("L51").Value = 7400
$ws.Range("M51").Value = -1516
$ws.Range("N51").Value = -8368

$ws.Range("H62").Value = 1650
$ws.Range("I62").Value = 1300
$ws.Range("J62").Value = 2000
$ws.Range("K62").Value = 1300
$ws.Range("L62").Value = 2000
$ws.Range("M62").Value = -676
$ws.Range("N62").Value = -3248

$ws.Range("H65").Value = 1650
$ws.Range("I65").Value = 1300
$ws.Range("J65").Value = 2000
$ws.Range("K65").Value = 6500
$ws.Range("L65").Value = 10000
$ws.Range("M65").Value = -3380
$ws.Range("N65").Value = -16240

$ws.Range("H132").Value = 1169746.9
$ws.Range("I132").Value = 2557.0435
$ws.Range("J132").Value = 2582660.8
$ws.Range("K132").Value = 7671.130500000001
$ws.Range("L132").Value = 7747982.399999999
$ws.Range("M132").Value = -5141.130500000001
$ws.Range("N132").Value = -7753042.399999999

$ws.Range("H137").Value = 2117.611
$ws.Range("I137").Value = 1361.1875
$ws.Range("J137").Value = 2722.75
$ws.Range("K137").Value = 4083.5625
$ws.Range("L137").Value = 8168.25
$ws.Range("M137").Value = -1533.5625
$ws.Range("N137").Value = -13268.25

$ws.Range("H138").Value = 5716002
$ws.Range("I138").Value = 1369.5
$ws.Range("J138").Value = 15386919
$ws.Range("K138").Value = 4108.5
$ws.Range("L138").Value = 46160757
$ws.Range("M138").Value = 1031.5
$ws.Range("N138").Value = -46171037

$ws.Range("H141").Value = 3071.1428
$ws.Range("I141").Value = 1916.3334
$ws.Range("K141").Value = 5749.0002
$ws.Range("M141").Value = -569.0002000000004

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()

$ws.Range("H32").Value = 20622.186
$ws.Range("I32").Value = 20073.54
$ws.Range("J32").Value = 22320.38
$ws.Range("K32").Value = 20073.54
$ws.Range("L32").Value = 22320.38
$ws.Range("M32").Value = -19786.54
$ws.Range("N32").Value = -22894.38

$ws.Range("H61").Value = 60440.676
$ws.Range("I61").Value = 34549.434
$ws.Range("K61").Value = 34549.434
$ws.Range("M61").Value = -34337.434

$ws.Range("H97").Value = 3677341.2
$ws.Range("I97").Value = 4465179.5
$ws.Range("J97").Value = 763.6667
$ws.Range("K97").Value = 4465179.5
$ws.Range("L97").Value = 763.6667
$ws.Range("M97").Value = -4464683.5
$ws.Range("N97").Value = -1755.6667

$ws.Range("H105").Value = 0
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

$ws.Range("H136").Value = 60440.676
$ws.Range("I136").Value = 34549.434
$ws.Range("K136").Value = 103648.302
$ws.Range("M136").Value = -101098.302

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("L25").Value = 0
$ws.Range("M25").ClearContents()
$ws.Range("N25").ClearContents()

$ws.Range("H58").Value = 2396.818
$ws.Range("I58").Value = 2396.818
$ws.Range("K58").Value = 2396.818
$ws.Range("M58").Value = -2193.818

$ws.Range("H136").Value = 2396.818
$ws.Range("I136").Value = 2396.818
$ws.Range("K136").Value = 7190.454000000001
$ws.Range("M136").Value = -4640.454000000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 555.6923
$ws.Range("I5").Value = 181.35294
$ws.Range("J5").Value = 1262.7778
$ws.Range("K5").Value = 544.05882
$ws.Range("L5").Value = 3788.3334
$ws.Range("M5").Value = -432.05882
$ws.Range("N5").Value = -4012.3334

$ws.Range("H10").Value = 2743.5833
$ws.Range("I10").Value = 115.875
$ws.Range("J10").Value = 7999
$ws.Range("K10").Value = 347.625
$ws.Range("L10").Value = 23997
$ws.Range("M10").Value = -208.625
$ws.Range("N10").Value = -24275

$ws.Range("H12").Value = 64.2
$ws.Range("J12").Value = 65.64706
$ws.Range("L12").Value = 196.94118
$ws.Range("N12").Value = -542.94118

$ws.Range("H61").Value = 1696.6666
$ws.Range("I61").Value = 90
$ws.Range("J61").Value = 2500
$ws.Range("K61").Value = 270
$ws.Range("L61").Value = 7500
$ws.Range("M61").Value = -55
$ws.Range("N61").Value = -7930

$ws.Range("H107").Value = 615.4872
$ws.Range("I107").Value = 512.65216
$ws.Range("J107").Value = 763.3125
$ws.Range("K107").Value = 1537.95648
$ws.Range("L107").Value = 2289.9375
$ws.Range("M107").Value = 382.0435200000002
$ws.Range("N107").Value = -6129.9375

$ws.Range("H122").Value = 786.4722
$ws.Range("I122").Value = 254.78572
$ws.Range("J122").Value = 1124.8182
$ws.Range("K122").Value = 2293.07148
$ws.Range("L122").Value = 10123.3638
$ws.Range("M122").Value = 156.9285199999999
$ws.Range("N122").Value = -15023.3638

$ws.Range("H131").Value = 1165.409
$ws.Range("I131").Value = 455.8
$ws.Range("J131").Value = 1374.1177
$ws.Range("K131").Value = 1367.4
$ws.Range("L131").Value = 4122.3531
$ws.Range("M131").Value = 3672.6
$ws.Range("N131").Value = -14202.3531

$ws.Range("H135").Value = 555.6923
$ws.Range("I135").Value = 181.35294
$ws.Range("J135").Value = 1262.7778
$ws.Range("K135").Value = 1632.17646
$ws.Range("L135").Value = 11365.0002
$ws.Range("M135").Value = 902.8235400000001
$ws.Range("N135").Value = -16435.0002

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 48954.652
$ws.Range("I70").Value = 129562.25
$ws.Range("J70").Value = 5963.933
$ws.Range("K70").Value = 129562.25
$ws.Range("L70").Value = 5963.933
$ws.Range("M70").Value = -129292.25
$ws.Range("N70").Value = -6503.933

$ws.Range("H73").Value = 48954.652
$ws.Range("I73").Value = 129562.25
$ws.Range("J73").Value = 5963.933
$ws.Range("K73").Value = 129562.25
$ws.Range("L73").Value = 5963.933
$ws.Range("M73").Value = -128626.25
$ws.Range("N73").Value = -7835.933

$ws.Range("H122").Value = 2554.8333
$ws.Range("I122").Value = 2193.75
$ws.Range("J122").Value = 3277
$ws.Range("K122").Value = 6581.25
$ws.Range("L122").Value = 9831
$ws.Range("M122").Value = -4131.25
$ws.Range("N122").Value = -14731

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 984.9
$ws.Range("I82").Value = 566.3333
$ws.Range("J82").Value = 1164.2858
$ws.Range("K82").Value = 566.3333
$ws.Range("L82").Value = 1164.2858
$ws.Range("M82").Value = -205.3333
$ws.Range("N82").Value = -1886.2858

$ws.Range("H85").Value = 984.9
$ws.Range("I85").Value = 566.3333
$ws.Range("J85").Value = 1164.2858
$ws.Range("K85").Value = 566.3333
$ws.Range("L85").Value = 1164.2858
$ws.Range("M85").Value = 681.6667
$ws.Range("N85").Value = -3660.2858

$ws.Range("H100").Value = 1527.6364
$ws.Range("I100").Value = 1234
$ws.Range("J100").Value = 1880
$ws.Range("K100").Value = 1234
$ws.Range("L100").Value = 1880
$ws.Range("M100").Value = -693
$ws.Range("N100").Value = -2962

$ws.Range("H132").Value = 105307.7
$ws.Range("I132").Value = 5699.5
$ws.Range("J132").Value = 130209.75
$ws.Range("K132").Value = 17098.5
$ws.Range("L132").Value = 390629.25
$ws.Range("M132").Value = -14568.5
$ws.Range("N132").Value = -395689.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
